$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.335.42'
$ws.Range("E2").Value = '  -3.06%  '
$ws.Range("D3").Value = '3.499.63'
$ws.Range("E3").Value = '  -4.66%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'603.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.14%  '
$ws.Range("D6").Value = "'149.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.21%  '
$ws.Range("D7").Value = '3.498.37'
$ws.Range("E7").Value = '  -4.63%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("D11").Value = "'6.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.54%  '
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("D14").Value = '4.091.99'
$ws.Range("E14").Value = '  -4.59%  '
$ws.Range("D15").Value = "'31.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.85%  '
$ws.Range("D16").Value = '3.503.10'
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("D17").Value = '67.273.66'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").Value = "'15.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.52%  '
$ws.Range("D21").Value = "'446.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.68%  '
$ws.Range("D22").Value = "'8.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -12.71%  '
$ws.Range("D23").Value = "'0.619"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.54%  '
$ws.Range("D24").Value = "'77.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.62%  '
$ws.Range("D25").Value = "'0.0000129"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.98%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '3.639.72'
$ws.Range("E27").Value = '  -4.65%  '
$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.48%  '
$ws.Range("D29").Value = "'8.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.22%  '
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("E31").Value = '  -6.89%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("E34").Value = '  -3.46%  '
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("D36").Value = '3.488.08'
$ws.Range("E36").Value = '  -4.96%  '
$ws.Range("E37").Value = '  -6.83%  '
$ws.Range("D38").Value = "'7.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = "'2.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").Value = "'173.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = "'5.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.57%  '
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("D46").Value = "'45.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("D47").Value = "'27.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("E48").Value = '  +6.05%  '
$ws.Range("E49").Value = '  -5.36%  '
$ws.Range("D50").Value = "'7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.15%  '
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.99%  '
